$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")
$ws.Activate()

# Update G9 (status of "Prototyping Study Report") to "Ready for printing"
$ws.Range("G9").Value = "Ready for printing"
$ws.Range("G9").WrapText = $false

# Reset row 9 height to default (remove explicit 63.75 height) since text no longer needs wrapping
$ws.Rows(9).AutoFit()

# Theme text color (Text 1 / dk1) darkened from #363636 to pure black #000000
$wb.Theme.ThemeColorScheme.Item(1).RGB = 0

# Update the active selection to C9
$ws.Range("C9").Select()
